# EIA Table 1.18.B monthly refresh: October 2016 YTD -> November 2016 YTD
# (and corresponding YTD-through-November data refresh for the rows whose
# underlying monthly figures changed between the October and November runs)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_1_18_B")

# --- Title / header text updates -------------------------------------------------

$ws.Range("A2").Value = "by State, by Sector, Year-to-Date through November 2016 and 2015 (Thousand Megawatthours)"

$ws.Range("B6").Value = "November 2016 YTD"
$ws.Range("E6").Value = "November 2016 YTD"
$ws.Range("G6").Value = "November 2016 YTD"
$ws.Range("I6").Value = "November 2016 YTD"
$ws.Range("K6").Value = "November 2016 YTD"

$ws.Range("C6").Value = "November 2015 YTD"
$ws.Range("F6").Value = "November 2015 YTD"
$ws.Range("H6").Value = "November 2015 YTD"
$ws.Range("J6").Value = "November 2015 YTD"
$ws.Range("L6").Value = "November 2015 YTD"

# --- Data updates ------------------------------------------------------------------

# Row 32: South Atlantic
$ws.Range("B32").Value = 73
$ws.Range("C32").Value = 104
$ws.Range("D32").Value = -0.302
$ws.Range("E32").Value = 73
$ws.Range("F32").Value = 104

# Row 35: Florida
$ws.Range("B35").Value = 73
$ws.Range("C35").Value = 104
$ws.Range("D35").Value = -0.302
$ws.Range("E35").Value = 73
$ws.Range("F35").Value = 104

# Row 52: Mountain
$ws.Range("B52").Value = 864
$ws.Range("C52").Value = 796
$ws.Range("D52").Value = 0.086
$ws.Range("G52").Value = 864
$ws.Range("H52").Value = 796

# Row 53: Arizona
$ws.Range("B53").Value = 622
$ws.Range("C53").Value = 689
$ws.Range("D53").Value = -0.098
$ws.Range("G53").Value = 622
$ws.Range("H53").Value = 689

# Row 57: Nevada
$ws.Range("B57").Value = 242
$ws.Range("C57").Value = 107
$ws.Range("D57").Value = 1.272
$ws.Range("G57").Value = 242
$ws.Range("H57").Value = 107

# Row 61: Pacific Contiguous
$ws.Range("B61").Value = 2356
$ws.Range("C61").Value = 2201
$ws.Range("D61").Value = 0.07
$ws.Range("G61").Value = 2356
$ws.Range("H61").Value = 2201

# Row 62: California
$ws.Range("B62").Value = 2356
$ws.Range("C62").Value = 2201
$ws.Range("D62").Value = 0.07
$ws.Range("G62").Value = 2356
$ws.Range("H62").Value = 2201

# Row 68: U.S. Total
$ws.Range("B68").Value = 3293
$ws.Range("C68").Value = 3101
$ws.Range("D68").Value = 0.062
$ws.Range("E68").Value = 73
$ws.Range("F68").Value = 104
$ws.Range("G68").Value = 3220
$ws.Range("H68").Value = 2997
